$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 148, shifting existing rows 148:171 down to 150:173
$ws.Rows.Item(148).Resize(2).Insert()

# Copy style (date number format) of column D from the row that is now 150 (originally 148)
$ws.Range("D148:D149").Value = $ws.Range("D150").Value
$ws.Range("D148:D149").NumberFormat = $ws.Range("D150").NumberFormat

# Fill new row 148 (Primera) with the new weekly data
$ws.Range("A148").Value = 8
$ws.Range("B148").Value = "Terminal La Palmera de La Serena"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44476
$ws.Range("E148").Value = 4
$ws.Range("F148").Value = 100114014
$ws.Range("G148").Value = "Betarraga"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 3100
$ws.Range("K148").Value = 450
$ws.Range("L148").Value = 500
$ws.Range("M148").Value = 475
$ws.Range("N148").Value = "`$/paquete 3 unidades"
$ws.Range("O148").Value = "Provincia del Elquí"
$ws.Range("P148").Value = 158
$ws.Range("Q148").Value = 3
$ws.Range("R148").Value = "Hortaliza"

# Fill new row 149 (Segunda) with the new weekly data
$ws.Range("A149").Value = 8
$ws.Range("B149").Value = "Terminal La Palmera de La Serena"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 44476
$ws.Range("E149").Value = 4
$ws.Range("F149").Value = 100114014
$ws.Range("G149").Value = "Betarraga"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Segunda"
$ws.Range("J149").Value = 1400
$ws.Range("K149").Value = 350
$ws.Range("L149").Value = 400
$ws.Range("M149").Value = 375
$ws.Range("N149").Value = "`$/paquete 3 unidades"
$ws.Range("O149").Value = "Provincia del Elquí"
$ws.Range("P149").Value = 125
$ws.Range("Q149").Value = 3
$ws.Range("R149").Value = "Hortaliza"
